# Generate Report for Handback
# Refreshes the handoff/handback timestamps for the "41261206..." row
# (row 3) in each per-language sheet, and mirrors the newest timestamp
# back onto the Overview sheet's "Latest HO Xliff Generate Date" column.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn: row 3 = 41261206-05c9-4c23-85d7-253336fa2ef9
# Column H = Correspond Handoff Datetime, Column K = Correspond Handback DateTime
$zhcn.Range("H3").Value = "2016-08-24 22:48:37"
$zhcn.Range("K3").Value = "2016-08-24 22:48:54"

# de-de: row 3 = 41261206-05c9-4c23-85d7-253336fa2ef9
$dede.Range("H3").Value = "2016-08-24 22:48:42"
$dede.Range("K3").Value = "2016-08-24 22:49:04"

# Overview: row 3 = 41261206-05c9-4c23-85d7-253336fa2ef9, column G = Latest HO Xliff Generate Date
$overview.Range("G3").Value = "2016-08-24 22:48:42"
